$d = $word.ActiveDocument

function Break-Join([string]$findText, [string]$replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

Break-Join 'e Polímeros.- Consideraç' 'e Polímeros.^l- Consideraç'
Break-Join 's avançados.- Polímeros ' 's avançados.^l- Polímeros '
Break-Join 'xcepcionais.- Polímeros ' 'xcepcionais.^l- Polímeros '
Break-Join 'ia de forma.- Aplicações' 'ia de forma.^l- Aplicações'
Break-Join 'olímeros em:- transporte' 'olímeros em:^l- transporte'
Break-Join 'g delivery);- transforma' 'g delivery);^l- transforma'
Break-Join 'ransfecção);- próteses d' 'ransfecção);^l- próteses d'
Break-Join ' em humanos;- nanorreato' ' em humanos;^l- nanorreato'
Break-Join 'es/catálise;- descontami' 'es/catálise;^l- descontami'
Break-Join 'io-ambiente;- eletrônica' 'io-ambiente;^l- eletrônica'
Break-Join ' condutores;- agricultur' ' condutores;^l- agricultur'
Break-Join ' agricultura- revestimen' ' agricultura^l- revestimen'
Break-Join 'superfícies.- recuperaçã' 'superfícies.^l- recuperaçã'
Break-Join 'de petróleo.- Polímeros ' 'de petróleo.^l- Polímeros '
Break-Join 'idegradáveis- Polímeros ' 'idegradáveis^l- Polímeros '
Break-Join 'iopolímeros.- Relações e' 'iopolímeros.^l- Relações e'
